$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D for the two newly reported
# quarters; this shifts the existing D:K data right to F:M.
$ws.Range("D:E").Insert()

# The inserted columns pick up a blank default style; copy the number
# formatting back from column F (the original column D, now shifted)
# so dates stay dates and figures stay figures.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns (D = newest quarter, E = prior quarter)
# with the newly reported figures for every data row.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 3344000
$ws.Range("E8").Value = 3390000
$ws.Range("D9").Value = 2534000
$ws.Range("E9").Value = 2534000
$ws.Range("D10").Value = 810000
$ws.Range("E10").Value = 856000
$ws.Range("D12").Value = 7000
$ws.Range("E12").Value = 7000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 138000
$ws.Range("E14").Value = 177000
$ws.Range("D15").Value = 174000
$ws.Range("E15").Value = 173000
$ws.Range("D17").Value = 2912000
$ws.Range("E17").Value = 2949000
$ws.Range("D18").Value = 432000
$ws.Range("E18").Value = 441000
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = -32000
$ws.Range("E20").Value = -2000
$ws.Range("D21").Value = 574000
$ws.Range("E21").Value = 611000
$ws.Range("D22").Value = 31000
$ws.Range("E22").Value = 33000
$ws.Range("D23").Value = 369000
$ws.Range("E23").Value = 406000
$ws.Range("D24").Value = 157000
$ws.Range("E24").Value = 251000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 212000
$ws.Range("E26").Value = 155000
$ws.Range("D27").Value = 43000
$ws.Range("E27").Value = -41000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 32000
$ws.Range("E32").Value = 2000
$ws.Range("D33").Value = 43000
$ws.Range("E33").Value = -41000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 43000
$ws.Range("E35").Value = -41000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 1113000
$ws.Range("E41").Value = 1022000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1003000
$ws.Range("E43").Value = 1193000
$ws.Range("D44").Value = 1644000
$ws.Range("E44").Value = 1666000
$ws.Range("D45").Value = 374000
$ws.Range("E45").Value = 312000
$ws.Range("D46").Value = 4134000
$ws.Range("E46").Value = 4193000
$ws.Range("D47").Value = 1360000
$ws.Range("E47").Value = 1381000
$ws.Range("D48").Value = 8327000
$ws.Range("E48").Value = 8355000
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 2117000
$ws.Range("E52").Value = 2256000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 15938000
$ws.Range("E54").Value = 16185000
$ws.Range("D57").Value = 1663000
$ws.Range("E57").Value = 1711000
$ws.Range("D58").Value = 1000
$ws.Range("E58").Value = 4000
$ws.Range("D59").Value = 1255000
$ws.Range("E59").Value = 1289000
$ws.Range("D60").Value = 2919000
$ws.Range("E60").Value = 3004000
$ws.Range("D61").Value = 1801000
$ws.Range("E61").Value = 1820000
$ws.Range("D62").Value = 3824000
$ws.Range("E62").Value = 4076000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 10549000
$ws.Range("E66").Value = 10969000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 341000
$ws.Range("E72").Value = 298000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 5389000
$ws.Range("E76").Value = 5216000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 43000
$ws.Range("E81").Value = -41000
$ws.Range("D83").Value = 174000
$ws.Range("E83").Value = 172000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 535000
$ws.Range("E89").Value = 288000
$ws.Range("D91").Value = -148000
$ws.Range("E91").Value = -82000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -148000
$ws.Range("E94").Value = -83000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -294000
$ws.Range("E100").Value = -280000
$ws.Range("D101").Value = -3000
$ws.Range("E101").Value = 6000
$ws.Range("D102").Value = 90000
$ws.Range("E102").Value = -69000

# Data correction: "Net Borrowings" values for the quarters now in
# columns I and J (previously G/H before the insert) were restated.
$ws.Range("I91").Value = -96000
$ws.Range("J91").Value = -88000

